$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for rows 2-6 as part of the
# repull/mean recalculation pass.
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 2
